$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (preserve exact source formatting,
# e.g. trailing zeros like "1.00") instead of being auto-converted to numbers.
$textCells = @('D5', 'D6', 'D8', 'D9', 'D14', 'D16', 'D18', 'D19', 'D22', 'D23', 'D25', 'D26', 'D27', 'D31', 'D34', 'D37', 'D39', 'D42', 'D43', 'D45', 'D46', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates per diff
$ws.Range('D2').Value = '26.776.42'
$ws.Range('E2').Value = '  -2.41%  '
$ws.Range('D3').Value = '1.561.29'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '206.07'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('D6').Value = '0.489'
$ws.Range('E6').Value = '  -1.98%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '21.95'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '0.248'
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').Value = '1.781.80'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '1.536.56'
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('D14').Value = '3.74'
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').Value = '61.55'
$ws.Range('E16').Value = '  -2.80%  '
$ws.Range('D17').Value = '26.758.80'
$ws.Range('E17').Value = '  -2.43%  '
$ws.Range('D18').Value = '214.91'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').Value = '7.34'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').Value = '0.0₃0677'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = '4.09'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').Value = '9.34'
$ws.Range('E23').Value = '  -2.22%  '
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('D25').Value = '152.55'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').Value = '6.75'
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('D27').Value = '14.88'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('D31').Value = '0.0462'
$ws.Range('E31').Value = '  -1.57%  '
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('D33').Value = '1.387.39'
$ws.Range('E33').Value = '  +1.97%  '
$ws.Range('D34').Value = '2.91'
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('D37').Value = '0.929'
$ws.Range('E37').Value = '  -4.42%  '
$ws.Range('E38').Value = '  -2.79%  '
$ws.Range('D39').Value = '0.808'
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('E40').Value = '  -3.23%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '0.996'
$ws.Range('E42').Value = '  +2.29%  '
$ws.Range('D43').Value = '5.42'
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('D45').Value = '63.31'
$ws.Range('E45').Value = '  -1.24%  '
$ws.Range('D46').Value = '1.76'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('D47').Value = '1.694.90'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').Value = '85.67'
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0493'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0945'
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.12%  '
